# Adds a new "migration" row (row 31) to the users-import worksheet so the
# evaluation history of a student migrating from min2 -> cin (etc.) is kept.
#
# New row layout (headers are A=firstname B=lastname C=email D=login
# E=roles F=groups G=period H=comment):
#   A31 migra                     (firstname)
#   B31 tor                       (lastname)
#   C31 migrator@eduvaud.ch       (email, mailto hyperlink)
#   D31 migrator@eduvaud.ch       (login, mailto hyperlink)
#   E31 eleve                     (roles)
#   F31 cin2a                     (groups)
#   G31 <empty, date-formatted>   (period)
#   H31 migration depuis min2     (comment)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in values in the same order the original workbook ended up storing
# its shared strings (comment, then email/login, then lastname, firstname).
$ws.Range("H31").Value = "migration depuis min2"
$ws.Range("C31").Value = "migrator@eduvaud.ch"
$ws.Range("D31").Value = "migrator@eduvaud.ch"
$ws.Range("B31").Value = "tor"
$ws.Range("A31").Value = "migra"
$ws.Range("E31").Value = "eleve"
$ws.Range("F31").Value = "cin2a"

# G31 stays empty but keeps the same short-date number format used by the
# rest of the "period" column (e.g. G29).
$ws.Range("G31").NumberFormat = "m/d/yy"

# Email/login columns are mailto hyperlinks, like the rest of the sheet.
$ws.Hyperlinks.Add($ws.Range("C31"), "mailto:migrator@eduvaud.ch")
$ws.Hyperlinks.Add($ws.Range("D31"), "mailto:migrator@eduvaud.ch")

# Hyperlinks.Add() changes the cell style to its own default; restore the
# workbook's usual "Lien hypertexte" style used by C/D columns elsewhere.
$ws.Range("C31").Style = "Lien hypertexte"
$ws.Range("D31").Style = "Lien hypertexte"

# Reflect the new selection/active cell left by the edit.
$ws.Range("A31:I31").Select()
